$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "1220194200662"
$ws.Range("F3").Select()
